$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.431.62"
$ws.Range("E2").Value = "'  +1.02%  "
$ws.Range("D3").Value = "'1.673.57"
$ws.Range("E3").Value = "'  +1.09%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("D5").Value = "'221.36"
$ws.Range("E5").Value = "'  +1.63%  "
$ws.Range("E6").Value = "'  +0.85%  "
$ws.Range("D7").Value = "'1.011"
$ws.Range("E7").Value = "'  +0.52%  "
$ws.Range("D8").Value = "'0.2669"
$ws.Range("E8").Value = "'  +1.81%  "
$ws.Range("D9").Value = "'0.06390"
$ws.Range("E9").Value = "'  +1.13%  "
$ws.Range("D10").Value = "'20.89"
$ws.Range("E10").Value = "'  +2.46%  "
$ws.Range("D11").Value = "'0.07858"
$ws.Range("E11").Value = "'  +0.58%  "
$ws.Range("D12").Value = "'4.536"
$ws.Range("E12").Value = "'  +0.45%  "
$ws.Range("D13").Value = "'1.682.59"
$ws.Range("E13").Value = "'  +2.32%  "
$ws.Range("D14").Value = "'1.903.52"
$ws.Range("E14").Value = "'  +1.10%  "
$ws.Range("D15").Value = "'0.5623"
$ws.Range("E15").Value = "'  +2.50%  "
$ws.Range("D16").Value = "'" + "0.0" + [char]0x2085 + "8205"
$ws.Range("E16").Value = "'  +0.64%  "
$ws.Range("D17").Value = "'66.25"
$ws.Range("E17").Value = "'  +1.41%  "
$ws.Range("D18").Value = "'26.451.02"
$ws.Range("E18").Value = "'  +1.19%  "
$ws.Range("E19").Value = "'  +0.56%  "
$ws.Range("D20").Value = "'4.731"
$ws.Range("E20").Value = "'  +2.93%  "
$ws.Range("D21").Value = "'197.57"
$ws.Range("E21").Value = "'  +3.47%  "
$ws.Range("D22").Value = "'10.32"
$ws.Range("E22").Value = "'  +2.52%  "
$ws.Range("D23").Value = "'6.085"
$ws.Range("E23").Value = "'  +1.51%  "
$ws.Range("E24").Value = "'  +0.43%  "
$ws.Range("D25").Value = "'146.29"
$ws.Range("E25").Value = "'  +0.56%  "
$ws.Range("D26").Value = "'0.1230"
$ws.Range("E26").Value = "'  +0.43%  "
$ws.Range("D27").Value = "'7.260"
$ws.Range("E27").Value = "'  +0.77%  "
$ws.Range("D28").Value = "'16.22"
$ws.Range("E28").Value = "'  +1.50%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "'  +2.52%  "
$ws.Range("D30").Value = "'0.05929"
$ws.Range("E30").Value = "'  +3.73%  "
$ws.Range("D31").Value = "'1.291"
$ws.Range("E31").Value = "'  +1.37%  "
$ws.Range("D32").Value = "'3.566"
$ws.Range("E32").Value = "'  +0.51%  "
$ws.Range("D33").Value = "'3.341"
$ws.Range("E33").Value = "'  +2.33%  "
$ws.Range("D34").Value = "'1.615"
$ws.Range("E34").Value = "'  +1.62%  "
$ws.Range("D35").Value = "'0.9698"
$ws.Range("E35").Value = "'  +2.28%  "
$ws.Range("D36").Value = "'2.840"
$ws.Range("E36").Value = "'  +1.15%  "
$ws.Range("D37").Value = "'2.439"
$ws.Range("E37").Value = "'  +0.76%  "
$ws.Range("D38").Value = "'0.5844"
$ws.Range("E38").Value = "'  +2.13%  "
$ws.Range("D39").Value = "'0.01617"
$ws.Range("E39").Value = "'  +0.62%  "
$ws.Range("D40").Value = "'1.080.83"
$ws.Range("E40").Value = "'  +4.18%  "
$ws.Range("D41").Value = "'5.931"
$ws.Range("E41").Value = "'  +2.38%  "
$ws.Range("D42").Value = "'0.8660"
$ws.Range("E42").Value = "'  +1.80%  "
$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "'  +0.52%  "
$ws.Range("D44").Value = "'103.31"
$ws.Range("E44").Value = "'  -0.56%  "
$ws.Range("D45").Value = "'1.813.24"
$ws.Range("E45").Value = "'  +0.99%  "
$ws.Range("D46").Value = "'58.69"
$ws.Range("E46").Value = "'  +3.53%  "
$ws.Range("D47").Value = "'" + "0.0" + [char]0x2088 + "107"
$ws.Range("E47").Value = "'  +2.06%  "
$ws.Range("D48").Value = "'1.016"
$ws.Range("E48").Value = "'  +0.93%  "
$ws.Range("D49").Value = "'0.4416"
$ws.Range("E49").Value = "'  +1.37%  "
$ws.Range("D50").Value = "'8.003"
$ws.Range("E50").Value = "'  +2.03%  "
$ws.Range("D51").Value = "'0.05162"
$ws.Range("E51").Value = "'  +0.19%  "
